# Update Chris Morris innings-by-innings batting activity (runs/balls/fours)
# to reflect the latest figures captured "till excel form".
#
# Source cells in this sheet store numeric-looking figures as TEXT, so each
# write forces a Text number format before assigning the value (otherwise
# Excel would auto-coerce the string into a real number), then restores the
# cell style to Normal so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Value
    )
    $rng = $ws.Range($Address)
    $rng.NumberFormat = "@"
    $rng.Value = $Value
    $rng.Style = "Normal"
}

# Row 2 (innings 1): runs 4 -> 0, fours 1 -> 0
Set-TextValue "C2" "0"
Set-TextValue "E2" "0"

# Row 3 (innings 2): runs 3 -> 2, balls 4 -> 5
Set-TextValue "C3" "2"
Set-TextValue "D3" "5"

# Row 4 (innings 3): runs 0 -> 3, balls 2 -> 4
Set-TextValue "C4" "3"
Set-TextValue "D4" "4"

# Row 5 (innings 4): runs 2 -> 4, balls 5 -> 2, fours 0 -> 1
Set-TextValue "C5" "4"
Set-TextValue "D5" "2"
Set-TextValue "E5" "1"
